# Swap the order of "System" and the email address in column G
# ("System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System")
# for every row where that exact value appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
